$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00144822592324403
$ws.Range("C2").Value = 0.136133236784938
$ws.Range("D2").Value = 0.139029688631427
$ws.Range("E2").Value = 0.0738595220854453
$ws.Range("F2").Value = 0.00144822592324403
$ws.Range("G2").Value = 0.00144822592324403
$ws.Range("H2").Value = 0.864590876176684
$ws.Range("I2").Value = 0.00941346850108617
$ws.Range("J2").Value = 0.000724112961622013
$ws.Range("K2").Value = 0.862418537291817
$ws.Range("M2").Value = 0.0238957277335264
$ws.Range("N2").Value = 0.00217233888486604
$ws.Range("O2").Value = 0.0123099203475742
$ws.Range("P2").Value = 0.826212889210717
$ws.Range("Q2").Value = 0.00724112961622013
$ws.Range("R2").Value = 0.00506879073135409
$ws.Range("S2").Value = 0.99637943519189
$ws.Range("T2").Value = 0.0108616944243302
$ws.Range("U2").Value = 0.423606082548878
$ws.Range("V2").Value = 0.0246198406951484
$ws.Range("W2").Value = 0.0992034757422158
$ws.Range("X2").Value = 0.00796524257784214
$ws.Range("B3").Value = 0.952208544532947
$ws.Range("C3").Value = 0.853005068790731
$ws.Range("D3").Value = 0.0166545981173063
$ws.Range("E3").Value = 0.880521361332368
$ws.Range("F3").Value = 0.982621288921072
$ws.Range("G3").Value = 0.00217233888486604
$ws.Range("H3").Value = 0.0166545981173063
$ws.Range("I3").Value = 0.00289645184648805
$ws.Range("J3").Value = 0.986241853729182
$ws.Range("K3").Value = 0.0658942795076032
$ws.Range("L3").Value = 0.0275162925416365
$ws.Range("M3").Value = 0.00362056480811007
$ws.Range("N3").Value = 0.163649529326575
$ws.Range("O3").Value = 0.980448950036206
$ws.Range("Q3").Value = 0.0383779869659667
$ws.Range("R3").Value = 0.994207096307024
$ws.Range("S3").Value = 0.00217233888486604
$ws.Range("T3").Value = 0.942795076031861
$ws.Range("U3").Value = 0.0238957277335264
$ws.Range("V3").Value = 0.0238957277335264
$ws.Range("W3").Value = 0.0057929036929761
$ws.Range("X3").Value = 0.00217233888486604
$ws.Range("B4").Value = 0.00434467776973208
$ws.Range("C4").Value = 0.00144822592324403
$ws.Range("D4").Value = 0.769007965242578
$ws.Range("E4").Value = 0.0340333091962346
$ws.Range("F4").Value = 0.00144822592324403
$ws.Range("G4").Value = 0.00217233888486604
$ws.Range("H4").Value = 0.1151339608979
$ws.Range("I4").Value = 0.984793627805938
$ws.Range("J4").Value = 0.0101375814627082
$ws.Range("K4").Value = 0.0629978276611151
$ws.Range("L4").Value = 0.0123099203475742
$ws.Range("M4").Value = 0.00362056480811007
$ws.Range("N4").Value = 0.0101375814627082
$ws.Range("O4").Value = 0.00651701665459812
$ws.Range("P4").Value = 0.164373642288197
$ws.Range("Q4").Value = 0.0231716147719044
$ws.Range("R4").Value = 0.000724112961622013
$ws.Range("S4").Value = 0.00144822592324403
$ws.Range("T4").Value = 0.000724112961622013
$ws.Range("U4").Value = 0.545981173062998
$ws.Range("V4").Value = 0.0246198406951484
$ws.Range("W4").Value = 0.881969587255612
$ws.Range("X4").Value = 0.98551774076756
$ws.Range("B5").Value = 0.0419985517740768
$ws.Range("C5").Value = 0.00941346850108617
$ws.Range("D5").Value = 0.0738595220854453
$ws.Range("E5").Value = 0.0115858073859522
$ws.Range("F5").Value = 0.0144822592324403
$ws.Range("G5").Value = 0.994207096307024
$ws.Range("H5").Value = 0.00289645184648805
$ws.Range("I5").Value = 0.00289645184648805
$ws.Range("J5").Value = 0.00217233888486604
$ws.Range("K5").Value = 0.00868935553946416
$ws.Range("L5").Value = 0.960173787110789
$ws.Range("M5").Value = 0.968863142650253
$ws.Range("N5").Value = 0.823316437364229
$ws.Range("O5").Value = 0.000724112961622013
$ws.Range("P5").Value = 0.00941346850108617
$ws.Range("Q5").Value = 0.931209268645909
$ws.Range("T5").Value = 0.0456191165821868
$ws.Range("U5").Value = 0.0057929036929761
$ws.Range("V5").Value = 0.926864590876177
$ws.Range("W5").Value = 0.0130340333091962
$ws.Range("X5").Value = 0.00362056480811007
